# Fix another issue in which the bootstrap capacitors were not connected
# to the source pins of the high-side FETs.
#
# The BOM's "100nF" capacitor group (and the groups it borrowed/lent
# references with) needs its reference list corrected, and the "Generated
# on" timestamp bumped to reflect the regenerated BOM.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference(s) column (F) for the capacitor groups that changed membership.
$ws.Range("F3").Value = "C1, C2, C4, C8, C20, C23, C30"
$ws.Range("F4").Value = "C3, C5, C6, C7, C9, C10, C12, C14, C15, C19"
$ws.Range("F5").Value = "C11, C16, C17, C18, C21, C22, C26, C28, C29, C31, C32, C33"
$ws.Range("F6").Value = "C13, C25, C34, C35"
$ws.Range("F7").Value = "C24"

# Regenerated-on date stamp (column B, summary block).
$ws.Range("B4").Value = "Generated on: 12/30/2020"

# Restore the last-used selection recorded by Excel on save.
$ws.Range("B16").Select()
